$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Test Case List")
$ws1.Range("C3").Value = ""
$ws1.Range("C3").Select()

$ws2 = $wb.Worksheets.Item("Test Data")
$ws2.Range("B3").Value = ""
$ws2.Range("B4").Value = ""
$ws2.Range("C3").Value = "FF"
$ws2.Range("C4").Value = "FF"
$ws2.Range("C3").Select()

$ws1.Activate()
$ws1.Range("C3").Select()
